$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 19:35"

# Row 4 - Estados Unidos: refreshed case counts
$ws.Range("B4").Value = 1602466
$ws.Range("C4").Value = 9743
$ws.Range("D4").Value = 371795
$ws.Range("E4").Value = 1135204
$ws.Range("G4").Value = 531
$ws.Range("H4").Value = 95467

# Row 33 - Irlanda: refreshed case counts
$ws.Range("B33").Value = 24391
$ws.Range("C33").Value = 76
$ws.Range("E33").Value = 1748
$ws.Range("G33").Value = 12
$ws.Range("H33").Value = 1583

# Row 41 - Israel: refreshed case counts
$ws.Range("B41").Value = 16683
$ws.Range("C41").Value = 16
$ws.Range("D41").Value = 13724
$ws.Range("E41").Value = 2680

# Row 105 - Sri Lanka: refreshed case counts
$ws.Range("B105").Value = 1047
$ws.Range("C105").Value = 19
$ws.Range("E105").Value = 434

# Rows 127/128 - Republica del Chad overtakes Sierra Leona in ranking
$ws.Range("A127").Value = "Republica del Chad"
$ws.Range("B127").Value = 588
$ws.Range("C127").Value = 23
$ws.Range("D127").Value = 186
$ws.Range("E127").Value = 344
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 58

$ws.Range("A128").Value = "Sierra Leona"
$ws.Range("B128").Value = 570
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 205
$ws.Range("E128").Value = 331
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 34

# Rows 199/200 - Belice overtakes Santa Lucia in ranking
$ws.Range("A199").Value = "Belice"
$ws.Range("B199").Value = 18
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 16
$ws.Range("E199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 2

$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 18
$ws.Range("E200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0

# Rows 209/210 - Groenlandia overtakes Seychelles in ranking (data tied)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Seychelles"

# Rows 214/215 - Sahara Occidental overtakes Bonaire, San Eustaquio y Saba (data tied)
$ws.Range("A214").Value = "Sahara Occidental"
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
